$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.475.57'
$ws.Range('E2').Value = '  +0.81%  '

$ws.Range('D3').Value = '3.607.28'
$ws.Range('E3').Value = '  +1.86%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '203.03'
$ws.Range('E5').Value = '  +4.26%  '

$ws.Range('E6').Value = '  -1.71%  '

$ws.Range('E7').Value = '  +0.59%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  +5.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.645'
$ws.Range('E10').Value = '  -0.38%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.61'
$ws.Range('E11').Value = '  +0.34%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000302'
$ws.Range('E12').Value = '  -0.61%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.70'
$ws.Range('E13').Value = '  +2.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '689.40'
$ws.Range('E14').Value = '  +16.32%  '

$ws.Range('D15').Value = '4.173.50'
$ws.Range('E15').Value = '  +1.86%  '

$ws.Range('D16').Value = '70.556.47'
$ws.Range('E16').Value = '  +0.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.10'
$ws.Range('E17').Value = '  -0.15%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.75'
$ws.Range('E18').Value = '  -0.36%  '

$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.592.32'
$ws.Range('E19').Value = '  +1.92%  '

$ws.Range('E20').Value = '  +0.50%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.997'
$ws.Range('E21').Value = '  +1.25%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.44'
$ws.Range('E22').Value = '  +3.56%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '110.05'
$ws.Range('E23').Value = '  +6.70%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.25'
$ws.Range('E24').Value = '  +2.16%  '

$ws.Range('E25').Value = '  -2.10%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.04'
$ws.Range('E26').Value = '  -0.55%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.60'
$ws.Range('E27').Value = '  -1.37%  '

$ws.Range('E28').Value = '  -0.69%  '

$ws.Range('E29').Value = '  +5.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.29'
$ws.Range('E30').Value = '  +3.33%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.49'
$ws.Range('E31').Value = '  +5.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.10'
$ws.Range('E32').Value = '  +0.71%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.33'
$ws.Range('E33').Value = '  -0.02%  '

$ws.Range('E34').Value = '  -1.22%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '63.81'
$ws.Range('E35').Value = '  +0.58%  '

$ws.Range('D36').Value = '0.0₃0850'
$ws.Range('E36').Value = '  +3.16%  '

$ws.Range('D37').Value = '3.858.25'
$ws.Range('E37').Value = '  +0.48%  '

$ws.Range('E38').Value = '  -0.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '510.08'
$ws.Range('E39').Value = '  -0.24%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.01'
$ws.Range('E40').Value = '  -6.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.60'
$ws.Range('E41').Value = '  +1.06%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.87'
$ws.Range('E42').Value = '  +1.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.384'
$ws.Range('E43').Value = '  -1.89%  '

$ws.Range('E44').Value = '  +3.44%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0472'
$ws.Range('E45').Value = '  +5.34%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.07'
$ws.Range('E46').Value = '  +9.14%  '

$ws.Range('E47').Value = '  +4.04%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.141'
$ws.Range('E48').Value = '  +0.93%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.64'
$ws.Range('E49').Value = '  +1.56%  '

$ws.Range('E50').Value = '  -0.25%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.84'
$ws.Range('E51').Value = '  +24.54%  '
